{"js": "// Update the date label and the 25 division problems (5 per row, 5 rows)\n// to the new values, in document order. Replacements are applied by\n// paragraph position (not by literal text search) because some of the\n// original values (e.g. \"71\u00f77=\") repeat but map to different new values.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst newTexts = [\n  \"2024-02-29 Thursday\",\n  \"25\u00f78=\", \"67\u00f77=\", \"44\u00f75=\", \"11\u00f72=\", \"41\u00f77=\",\n  \"78\u00f73=\", \"69\u00f77=\", \"28\u00f79=\", \"60\u00f76=\", \"27\u00f75=\",\n  \"97\u00f73=\", \"39\u00f79=\", \"86\u00f74=\", \"34\u00f72=\", \"36\u00f79=\",\n  \"58\u00f76=\", \"55\u00f78=\", \"40\u00f75=\", \"27\u00f77=\", \"21\u00f79=\",\n  \"77\u00f74=\", \"81\u00f74=\", \"85\u00f73=\", \"91\u00f77=\", \"44\u00f77=\"\n];\n\nconst items = paragraphs.items;\nlet idx = 0;\nfor (let i = 0; i < items.length && idx < newTexts.length; i++) {\n  const para = items[i];\n  para.load(\"text\");\n  await context.sync();\n  if (para.text !== \"\") {\n    para.insertText(newTexts[idx], Word.InsertLocation.replace);\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date label and the 25 division problems (5 per row, 5 rows)\n# to the new values, in document order. Replacements are applied by\n# paragraph position (not by literal text search) because some of the\n# original values (e.g. \"71\u00f77=\") repeat but map to different new values.\n$d = $word.ActiveDocument\n\n$newTexts = @(\n  \"2024-02-29 Thursday\",\n  \"25\u00f78=\", \"67\u00f77=\", \"44\u00f75=\", \"11\u00f72=\", \"41\u00f77=\",\n  \"78\u00f73=\", \"69\u00f77=\", \"28\u00f79=\", \"60\u00f76=\", \"27\u00f75=\",\n  \"97\u00f73=\", \"39\u00f79=\", \"86\u00f74=\", \"34\u00f72=\", \"36\u00f79=\",\n  \"58\u00f76=\", \"55\u00f78=\", \"40\u00f75=\", \"27\u00f77=\", \"21\u00f79=\",\n  \"77\u00f74=\", \"81\u00f74=\", \"85\u00f73=\", \"91\u00f77=\", \"44\u00f77=\"\n)\n\n$idx = 0\nforeach ($p in $d.Paragraphs) {\n    if ($idx -ge $newTexts.Count) { break }\n    $t = $p.Range.Text -replace \"[\\r\\x07]\", \"\"\n    if ($t -ne \"\") {\n        $p.Range.Text = $newTexts[$idx]\n        $idx++\n    }\n}\n"}
